# Week 15 logged + Week 16 simulated: update play-log strings and summary totals
$wb = $excel.ActiveWorkbook

# --- YDS sheet: append this period's per-play yardage logs ---
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = '4 4 27 -3 2 3 6 -1 8 3 5 6 2 6 6 6 6 3 0 -2 0 6 4 5 2 5 9 1 6 2 5 0 5 1 10 2 2 2 2 7 14 3 1 3 1 11 3 8 2 4 2 7 -1 4 1 -1 11 4 -4 7 3 1 7 4 7 6 0 1 0 5 7 4 6 5 3 6 2 8 6 6 -2 20 2 5 6 8 10 2 1 5 1 3 10 -1 1 -1 5 8 3 0 1 8 5 1 -7 13 5 3 3 1 4 1 2 3 3 2 5 6 16 1 9 2 11 2 -2 -2 16 2 1 4 8 1 3 5 3 3 2 2 15 1 7 8 6 4 4 11 2 5 3 8 5 4 5 2 2 6 1 4 38 2 6 4 5 3 0 2 6 27 12 5 2 8 7 13 0 2 10 4 3 2 2 5 3 6 5 0 4 4 7 0 4 5 13 12 12 4 4 1 3 1 5 -1 -8 10 5 2 3 15 4 4 -5 -1 -3 12 9 1 1 1 2 3 4 7 3 2 4 2 -1 5 0 3 4 3 10 2 -2 4 -3 15 5 2 11 1 8 3 4 0 1 4 5 11 20 2 0 11 9 4 9 4 2 8 -2 0 14 3 11 0 0 7 6 7 3 8 4 5 1 0 4 3 4 4 3 1 1 0 1 5 7 3 4 4 3 10 10 2 1 6 3 0 0 7 4 1 7 5 17 1 2 0 2 2 22 1 3 1 16 2 6 2 2 4 1 9 5 -4 -3 13 4 1 0 -6 7 2 3 8 1 6 0 4 3 2 0 1 3 3 11 2 4 6 5 0 13 6 5 12 8 4 6 7 2 0 0 2 6 5 2 4 1 3 3 5 4 -1 -1 8 -6 2 30 0 1 0 13'
$wsYDS.Range("B3").Value = '2 6 5 12 7 11 9 2 34 11 17 17 -2 4 7 5 22 12 4 4 16 11 8 12 0 24 12 7 6 10 7 4 4 8 18 9 0 3 10 7 56 19 2 7 9 5 7 6 0 11 17 10 73 14 7 4 15 1 11 15 37 4 1 6 9 3 19 -1 21 33 16 20 11 2 -2 11 23 10 17 24 19 4 30 20 9 18 8 7 17 13 7 16 10 7 -7 9 3 6 4 6 9 26 27 9 13 6 3 7 12 2 7 3 5 6 6 6 4 16 39 10 14 4 -4 -4 9 12 4 -2 7 6 12 40 0 5 5 14 4 17 28 9 12 5 7 10 3 11 2 19 10 3 12 13 6 15 8 27 17 2 13 -2 20 4 12 1 17 14 7 5 14 4 23 20 19 5 22 14 3 1 8 14 13 4 8 11 1 20 22 7 9 8 14 1 15 6 18 21 9 6 16 6 5 6 6 30 7 39 6 4 5 12 12 7 13 4 15 10 18 5 4 15 2 8 6 0 1 1 12 9 8 14 9 16 10 20 2 9 11 11 6 3 12 6 8 23 -1 13 6 7 4 7 6 5 6 5 6 0 10 35 10 3 4 10 9 10 12 6 14 5 12 2 3 19 9 43 -1 6 10 8 19 19 7 19 7 6 46 5 1 7 5 9 7 5 12 29 5 6 1 13 10 3 -2'
$wsYDS.Range("C2").Value = '11 13 3 3 1 7 -3 1 4 3 2 0 7 2 2 0 14 5 2 3 0 4 1 2 9 -3 5 2 15 -2 2 6 41 6 2 7 5 3 0 9 -1 -6 2 -9 7 2 2 5 11 1 1 2 0 4 6 8 1 -2 4 4 8 2 3 4 3 2 2 0 3 0 7 6 5 1 6 9 1 2 5 2 2 8 2 3 2 4 -3 10 2 5 -3 1 2 -3 3 3 13 5 2 0 13 1 2 -3 4 17 2 1 5 3 1 23 2 10 5 3 12 3 3 0 11 2 1 8 3 -2 1 7 2 0 6 3 3 17 3 6 1 6 4 0 5 11 4 3 3 3 10 5 2 -1 1 3 0 3 9 15 2 2 13 2 1 5 3 2 0 2 0 -1 8 4 2 6 2 2 7 1 1 7 2 -1 5 3 0 11 14 1 7 1 -2 5 6 0 3 5 8 2 1 3 4 8 6 8 6 3 0 -3 4 3 8 3 4 9 6 24 2 4 0 3 6 10 -1 2 7 5 0 0 1 3 0 1 12 4 3 7 1 3 0 17 7 3 15 -1 6 2 1 1 10 1 4 3 0 4 0 5 5 4 0 2 2 7 4 5 4 11 1 4 10 10 6 1 3 2 3 3 4 3 4 1 5 5 4 5 4 4 4 6 5 8 1 4 19 0 9 4 4 10 3 1 8 8 38 1 9 12 7 -3 4 3 15 4 2 3 4 4 0 8 4 -3 1'
$wsYDS.Range("C3").Value = '10 6 9 11 -3 7 2 11 4 18 16 10 22 1 6 17 4 17 11 7 11 23 6 19 17 3 17 7 19 20 9 12 16 13 4 4 8 2 6 17 3 4 2 8 12 19 33 13 21 10 8 21 13 23 11 13 1 28 3 7 9 7 17 4 7 14 18 2 41 7 13 22 11 5 11 9 5 6 13 10 11 6 5 19 4 13 15 7 15 10 5 12 42 8 8 17 6 12 15 5 7 14 14 6 7 8 7 5 21 4 11 72 20 11 10 12 49 14 -4 11 12 6 32 19 6 20 9 19 27 49 3 5 7 27 19 17 14 9 10 5 13 9 5 12 13 2 9 0 3 11 5 21 19 6 12 6 24 15 3 4 12 5 0 12 1 10 3 15 11 7 13 6 19 5 10 2 20 4 10 26 4 25 17 1 31 5 31 12 19 19 15 7 3 4 -4 9 2 11 14 9 18 7 2 3 29 9 8 3 6 7 11 4 5 14 15 8 6 2 15 4 0 4 3 40 8 13 10 5 3 11 2 19 7 9 4 10 16 27 7 6 4 8 11 7 2 55 6 13 9 39 13 8 2 7 6 1 12 1 4 13 13 4 7 32 1 8 9 3 10 5 3 9 7 5 0 18 10 -5 5 2 15 9 34 7 3 28 14 13 16 2 11 7 4 18 7 15 3 12 24 7 4 0 12 13 9 14 20 12 5 11 5 10 2 4 5 -1 19 28 16 45 6 5 8 17 21 20 5 9 7 9 7 34 19 17'

# --- ST sheet: append this period's per-play logs ---
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B4").Value = '67 63 61 41 66 65 70 68 73 58 54 66 61 64 62 62 60 64 57 58 65 59 62'
$wsST.Range("B5").Value = '19 23 20 0 24 32 22 27 23 17 19 34 21 27 27 27 19 24 21 15 31 23 17'
$wsST.Range("B6").Value = '16 31 18 21 27 16 23 21 0 19 28 23 1 22 26 15 23 19 30 24 25 23 28 19 23 30 12 31 28 20 47 28'
$wsST.Range("D3").Value = '56 45 37 55 48 50 47 52 42 55 38 48 39 45 56 51 47 66 60 55 51 66 44 45 56 49 34 48 37 53 47 33 54 47 43 53 59 57 55 53 42 41 49 36 60 44 46 54'
$wsST.Range("D4").Value = '0 3 0 11 8 0 0 6 33 0 0 5 2 4 13 14 2 0 31 25 0 15 10 4 10 0 0 0 0 0 0 0 6 0 0 0 0 4 0 1 0 0 0 0 11 0 1 10'
$wsST.Range("D5").Value = '14 8 0 0 0 0 0 0 6 0 0 0 0 0 10 0 0 0 0 0 0 12 0 0 10 0 0 0 0 0 0 0 0 0 0 0 8 0 0 16 0 13 0 0 0 0 0 0 0 0'

# --- Updated season-to-date summary totals ---

$wsOFF = $wb.Worksheets.Item("OFF")

# OFF
$wsOFF.Range("B2").Value = 7
$wsOFF.Range("C2").Value = 169
$wsOFF.Range("D2").Value = 9
$wsOFF.Range("F2").Value = 69
$wsOFF.Range("G2").Value = 50
$wsOFF.Range("J2").Value = 38
$wsOFF.Range("N2").Value = 17
$wsOFF.Range("O2").Value = 28
$wsOFF.Range("C3").Value = 167
$wsOFF.Range("E3").Value = 31
$wsOFF.Range("F3").Value = 76
$wsOFF.Range("G3").Value = 29
$wsOFF.Range("H3").Value = 24
$wsOFF.Range("I3").Value = 46
$wsOFF.Range("J3").Value = 48
$wsOFF.Range("L3").Value = 215
$wsOFF.Range("M3").Value = 145
$wsOFF.Range("Q3").Value = 431

$wsDEF = $wb.Worksheets.Item("DEF")

# DEF
$wsDEF.Range("B2").Value = 5
$wsDEF.Range("C2").Value = 156
$wsDEF.Range("D2").Value = 7
$wsDEF.Range("E2").Value = 7
$wsDEF.Range("F2").Value = 58
$wsDEF.Range("G2").Value = 33
$wsDEF.Range("I2").Value = 4
$wsDEF.Range("J2").Value = 26
$wsDEF.Range("N2").Value = 15
$wsDEF.Range("O2").Value = 13
$wsDEF.Range("P2").Value = 7
$wsDEF.Range("C3").Value = 168
$wsDEF.Range("D3").Value = 1
$wsDEF.Range("E3").Value = 25
$wsDEF.Range("F3").Value = 102
$wsDEF.Range("G3").Value = 32
$wsDEF.Range("H3").Value = 19
$wsDEF.Range("I3").Value = 45
$wsDEF.Range("J3").Value = 65
$wsDEF.Range("L3").Value = 238
$wsDEF.Range("M3").Value = 172
$wsDEF.Range("Q3").Value = 426

# ST
$wsST.Range("B2").Value = 64
$wsST.Range("D2").Value = 48
$wsST.Range("F2").Value = 96
$wsST.Range("G2").Value = 84
$wsST.Range("J2").Value = 45
$wsST.Range("K2").Value = 43
$wsST.Range("B3").Value = 41

$wsTURNS = $wb.Worksheets.Item("TURNS")

# TURNS
$wsTURNS.Range("C3").Value = 1
$wsTURNS.Range("D3").Value = 11
$wsTURNS.Range("E3").Value = 8
